$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A76").Value = 45907
$ws.Range("B76").Value = "四方坪站"
$ws.Range("C76").Value = 9789.4
$ws.Range("D76").Value = 8167.69
$ws.Range("E76").Value = 3313.4
$ws.Range("F76").Value = 400

$ws.Range("A77").Value = 45907
$ws.Range("B77").Value = "高岭站"
$ws.Range("C77").Value = 4548.07
$ws.Range("D77").Value = 3687.1
$ws.Range("E77").Value = 1075.24
$ws.Range("F77").Value = 158

$ws.Range("H72").Select()
